# Auto-update: append website scan results (rows 107-192) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(107, 1).Value = "https://xelere.com"
$ws.Cells.Item(107, 2).Value = "info@xelere.com"
$ws.Cells.Item(108, 1).Value = "https://iag.com.ar"
$ws.Cells.Item(108, 2).Value = "info@iag.com.ar"
$ws.Cells.Item(109, 1).Value = "https://keepcon.com"
$ws.Cells.Item(109, 2).Value = "info@keepcon.com"
$ws.Cells.Item(110, 1).Value = "https://trendingenieria.com.ar"
$ws.Cells.Item(110, 2).Value = "info@trendingenieria.com.ar"
$ws.Cells.Item(111, 1).Value = "https://imagecampus.com.ar"
$ws.Cells.Item(111, 2).Value = "info@imagecampus.edu.ar"
$ws.Cells.Item(112, 1).Value = "https://fls.org.ar"
$ws.Cells.Item(112, 2).Value = "info@fls.org.ar"
$ws.Cells.Item(113, 1).Value = "https://powersolution.com.ar"
$ws.Cells.Item(113, 2).Value = "info@powersolution.es, info@powersolution.com.ar, info@ps-iberia.net"
$ws.Cells.Item(114, 1).Value = "https://q4tech.com"
$ws.Cells.Item(114, 2).Value = "info@q4tech.com"
$ws.Cells.Item(115, 1).Value = "https://fwtv.tv"
$ws.Cells.Item(115, 2).Value = "info@fwtv.tv"
$ws.Cells.Item(116, 1).Value = "https://abeceb.com"
$ws.Cells.Item(116, 2).Value = "info@abeceb.com"
$ws.Cells.Item(117, 1).Value = "https://toribioachaval.com"
$ws.Cells.Item(117, 2).Value = "info@toribioachaval.com"
$ws.Cells.Item(118, 1).Value = "https://viditec.com"
$ws.Cells.Item(118, 2).Value = "marketing@viditec.com.ar, info@viditec.com.ar"
$ws.Cells.Item(119, 1).Value = "https://afluenta.com"
$ws.Cells.Item(119, 2).Value = "info@afluenta.com"
$ws.Cells.Item(120, 1).Value = "https://renova.com.ar"
$ws.Cells.Item(120, 2).Value = "info@renova.com.ar"
$ws.Cells.Item(121, 1).Value = "https://kornpropiedades.com.ar"
$ws.Cells.Item(121, 2).Value = "info@kornpropiedades.com.ar"
$ws.Cells.Item(122, 1).Value = "https://ucel.edu.ar"
$ws.Cells.Item(122, 2).Value = "info@ucel.edu.arHorario"
$ws.Cells.Item(123, 1).Value = "https://bodegadelfindelmundo.com"
$ws.Cells.Item(123, 2).Value = "info@bdfm.com.ar"
$ws.Cells.Item(124, 1).Value = "https://maprimed.com"
$ws.Cells.Item(124, 2).Value = "info@maprimed.com.ar"
$ws.Cells.Item(125, 1).Value = "https://sifeme.com"
$ws.Cells.Item(125, 2).Value = "info@sifemesa.com.ar"
$ws.Cells.Item(126, 1).Value = "https://siscard.com"
$ws.Cells.Item(126, 2).Value = "info@siscard.com.ar, info@siscardperu.pe"
$ws.Cells.Item(127, 1).Value = "https://gestiva.com.ar"
$ws.Cells.Item(127, 2).Value = "info@gestiva.com.arg, info@gestiva.com.ar"
$ws.Cells.Item(128, 1).Value = "https://obrasysistemas.com.ar"
$ws.Cells.Item(128, 2).Value = "info@obrasysistemas.com.ar"
$ws.Cells.Item(129, 1).Value = "https://loyal-solutions.com"
$ws.Cells.Item(129, 2).Value = "info@loyal-solutions.com"
$ws.Cells.Item(130, 1).Value = "https://grupolpa.com"
$ws.Cells.Item(130, 2).Value = "info@grupolpa.com"
$ws.Cells.Item(131, 1).Value = "https://gleba.com.ar"
$ws.Cells.Item(131, 2).Value = "info@gleba.com.ar"
$ws.Cells.Item(132, 1).Value = "https://cytcomunicaciones.com"
$ws.Cells.Item(132, 2).Value = "info@cytcomunicaciones.com.ar"
$ws.Cells.Item(133, 1).Value = "https://softland.com.ar"
$ws.Cells.Item(133, 2).Value = "info@softland.com.ar"
$ws.Cells.Item(134, 1).Value = "https://ieserh.com.ar"
$ws.Cells.Item(134, 2).Value = "info@ieserh.edu.ar"
$ws.Cells.Item(135, 1).Value = "https://tracegroup.com.ar"
$ws.Cells.Item(135, 2).Value = "info@tracegroup.com.ar"
$ws.Cells.Item(136, 1).Value = "https://cimientos.org"
$ws.Cells.Item(136, 2).Value = "info@cimientos.org"
$ws.Cells.Item(137, 1).Value = "https://laninia.com"
$ws.Cells.Item(137, 2).Value = "hello@laninia.com, jobs@laninia.com"
$ws.Cells.Item(138, 1).Value = "https://tucumanturismo.gob.ar"
$ws.Cells.Item(138, 2).Value = "info@riodearena.com, info@fincaalbarossa.com, info@bodegalasarcas.com.ar"
$ws.Cells.Item(139, 1).Value = "https://moellerip.com"
$ws.Cells.Item(139, 2).Value = "Marketing@moellerip.com, hello@moellerip.com"
$ws.Cells.Item(140, 1).Value = "https://dediego.com.ar"
$ws.Cells.Item(140, 2).Value = "info@dediego.com.ar"
$ws.Cells.Item(141, 1).Value = "https://brons.com.ar"
$ws.Cells.Item(141, 2).Value = "info@brons.com.ar"
$ws.Cells.Item(142, 1).Value = "https://enerminds.com"
$ws.Cells.Item(142, 2).Value = "info@prometium.com"
$ws.Cells.Item(143, 1).Value = "https://navicon.com.ar"
$ws.Cells.Item(143, 2).Value = "info@navicon.com.ar"
$ws.Cells.Item(144, 1).Value = "https://close-upinternational.com"
$ws.Cells.Item(144, 2).Value = "info@close-upinternational.com, info@closeupus.com"
$ws.Cells.Item(145, 1).Value = "https://bhp-global.com"
$ws.Cells.Item(145, 2).Value = "info@bhp-global.com"
$ws.Cells.Item(146, 1).Value = "https://fscnet.com.ar"
$ws.Cells.Item(146, 2).Value = "info@fsc.com.ar, info@fscnet.com.ar"
$ws.Cells.Item(147, 1).Value = "https://gruposanmiguel.com.ar"
$ws.Cells.Item(147, 2).Value = "info@gruposanmiguel.com.ar"
$ws.Cells.Item(148, 1).Value = "https://globalprocessing.com.ar"
$ws.Cells.Item(148, 2).Value = "info@globalprocessing.com.ar"
$ws.Cells.Item(149, 1).Value = "https://aadesa.com.ar"
$ws.Cells.Item(149, 2).Value = "info@aadesa.com.ar"
$ws.Cells.Item(150, 1).Value = "https://iquall.net"
$ws.Cells.Item(150, 2).Value = "hello@iquall.net"
$ws.Cells.Item(151, 1).Value = "https://geoagro.com"
$ws.Cells.Item(151, 2).Value = "info@geoagro.com"
$ws.Cells.Item(152, 1).Value = "https://alea.com.ar"
$ws.Cells.Item(152, 2).Value = "info@alea.com.ar, info@grupogaman.com.ar"
$ws.Cells.Item(153, 1).Value = "https://southend.com.ar"
$ws.Cells.Item(153, 2).Value = "info@southendcorp.com"
$ws.Cells.Item(154, 1).Value = "https://socmer.com.ar"
$ws.Cells.Item(154, 2).Value = "info@socmer.com.ar"
$ws.Cells.Item(155, 1).Value = "https://publicatulibro.com.ar"
$ws.Cells.Item(155, 2).Value = "info@publicatulibro.com.ar"
$ws.Cells.Item(156, 1).Value = "https://siemprearg.com"
$ws.Cells.Item(156, 2).Value = "info@siemprearg.com"
$ws.Cells.Item(157, 1).Value = "https://esama.com"
$ws.Cells.Item(157, 2).Value = "office@yoursite.com"
$ws.Cells.Item(158, 1).Value = "https://estudios-electricos.com"
$ws.Cells.Item(158, 2).Value = "info@estudios-electricos.com"
$ws.Cells.Item(159, 1).Value = "https://invgate.com"
$ws.Cells.Item(159, 2).Value = "info@invgate.com"
$ws.Cells.Item(160, 1).Value = "https://ott.edu.ar"
$ws.Cells.Item(160, 2).Value = "info@ott.edu.ar"
$ws.Cells.Item(161, 1).Value = "https://venturi.com.ar"
$ws.Cells.Item(161, 2).Value = "info@venturi.com.ar"
$ws.Cells.Item(162, 1).Value = "https://consensusgroup.net"
$ws.Cells.Item(162, 2).Value = "sales@sokosolutions.com, info@cstechlab.com"
$ws.Cells.Item(163, 1).Value = "https://humanagency.com.ar"
$ws.Cells.Item(163, 2).Value = "info@humanagency.com.ar"
$ws.Cells.Item(164, 1).Value = "https://telextorage.com"
$ws.Cells.Item(164, 2).Value = "info@telextorage.com"
$ws.Cells.Item(165, 1).Value = "https://winandwinnow.com"
$ws.Cells.Item(165, 2).Value = "contact@winandwinnow.com"
$ws.Cells.Item(166, 1).Value = "https://borealtech.com"
$ws.Cells.Item(166, 2).Value = "info@borealtech.com"
$ws.Cells.Item(167, 1).Value = "https://zbv.com.ar"
$ws.Cells.Item(167, 2).Value = "info@zbv.com.ar"
$ws.Cells.Item(168, 1).Value = "https://algeiba.com"
$ws.Cells.Item(168, 2).Value = "info@algeiba.com"
$ws.Cells.Item(169, 1).Value = "https://andesmarcargas.com"
$ws.Cells.Item(169, 2).Value = "info@andesmarcargas.com"
$ws.Cells.Item(170, 1).Value = "https://tb.com.ar"
$ws.Cells.Item(170, 2).Value = "Contact@tbgroup.tech, Info@tbgroup.tech"
$ws.Cells.Item(171, 1).Value = "https://properati.com"
$ws.Cells.Item(171, 2).Value = "info@properati.com"
$ws.Cells.Item(172, 1).Value = "https://dpisa.com.ar"
$ws.Cells.Item(172, 2).Value = "info@dpisa.com.ar"
$ws.Cells.Item(173, 1).Value = "https://celulosa.com.ar"
$ws.Cells.Item(173, 2).Value = "info@celulosa.com.ar"
$ws.Cells.Item(174, 1).Value = "https://sistemasactivos.com"
$ws.Cells.Item(174, 2).Value = "info@misitio.com, info@sistemasactivos.com"
$ws.Cells.Item(175, 1).Value = "https://bas.com.ar"
$ws.Cells.Item(175, 2).Value = "info@bas.com.ar"
$ws.Cells.Item(176, 1).Value = "https://silmag.com.ar"
$ws.Cells.Item(176, 2).Value = "info@silmag.com.ar"
$ws.Cells.Item(177, 1).Value = "https://tradelog.com.ar"
$ws.Cells.Item(177, 2).Value = "info@manderly.net, info@tradelog.com.ar"
$ws.Cells.Item(178, 1).Value = "https://tzedaka.org.ar"
$ws.Cells.Item(178, 2).Value = "info@tzedaka.org.ar"
$ws.Cells.Item(179, 1).Value = "https://toyotacfa.com.ar"
$ws.Cells.Item(179, 2).Value = "Info@otero-fraga.com.ar"
$ws.Cells.Item(180, 1).Value = "https://tinsa.com.ar"
$ws.Cells.Item(180, 2).Value = "info@tinsa.com.ar"
$ws.Cells.Item(181, 1).Value = "https://zillesrl.com.ar"
$ws.Cells.Item(181, 2).Value = "info@zoxisa.com.ar"
$ws.Cells.Item(182, 1).Value = "https://cimet.com"
$ws.Cells.Item(182, 2).Value = "info@cimet.com"
$ws.Cells.Item(183, 1).Value = "https://holistor.com.ar"
$ws.Cells.Item(183, 2).Value = "info@holistor.com.ar"
$ws.Cells.Item(184, 1).Value = "https://eficienciaempresaria.com"
$ws.Cells.Item(184, 2).Value = "info@eficienciaempresaria.com"
$ws.Cells.Item(185, 1).Value = "https://samsistemas.com.ar"
$ws.Cells.Item(185, 2).Value = "info@samsistemas.com.ar"
$ws.Cells.Item(186, 1).Value = "https://txtinternational.com"
$ws.Cells.Item(186, 2).Value = "info@txtinternational.com"
$ws.Cells.Item(187, 1).Value = "https://grupofarallon.com"
$ws.Cells.Item(187, 2).Value = "info@grupofarallon.com"
$ws.Cells.Item(188, 1).Value = "https://cetecsudamericana.com.ar"
$ws.Cells.Item(188, 2).Value = "info@penielsrl.com.ar, info@dsi.com.ar"
$ws.Cells.Item(189, 1).Value = "https://wineem.com.ar"
$ws.Cells.Item(189, 2).Value = "info@wineem.com.ar"
$ws.Cells.Item(190, 1).Value = "https://ashira-sa.com.ar"
$ws.Cells.Item(190, 2).Value = "info@amym-com-ar.png, info@amym.com.ar, info@ashira-sa.com.ar"
$ws.Cells.Item(191, 1).Value = "https://sg-sa.com.ar"
$ws.Cells.Item(191, 2).Value = "info@simpliafacility.com"
$ws.Cells.Item(192, 1).Value = "https://endeavor.org.ar"
$ws.Cells.Item(192, 2).Value = "support@everlytic.com, info@endeavor.org.ar"
